# Updated cryptos list on Sun Nov 19 20:11:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh price / volume(1h) figures for existing rows ---
$ws.Range("D2").Value = "37.032.60"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.986.60"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "245.95"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").Value = "61.53"
$ws.Range("E7").Value = "  +3.84%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "0.0806"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "14.95"
$ws.Range("E12").Value = "  +8.90%  "
$ws.Range("D13").Value = "22.41"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "2.276.87"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").Value = "1.990.18"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "36.950.15"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "70.44"
$ws.Range("D20").Value = "0.0₃0866"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "5.19"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "230.57"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "9.33"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").Value = "164.13"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "19.61"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +18.49%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "4.55"
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "5.57"
$ws.Range("E39").Value = "  -7.64%  "
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "16.48"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "1.379.68"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "90.40"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "7.31"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "2.05"
$ws.Range("E49").Value = "  +15.53%  "

# --- Rank #48/#49 swapped places: MXToken now ranks above MultiversX ---
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "46.44"
$ws.Range("E51").Value = "  +5.69%  "
